# Atualização de bases das ligas, do dia: 04-04-2024 às 23:22
#
# This script:
#  1. Swaps the (B:AC) data between several pairs of rows in the single
#     worksheet (the "id"/A column stays put, everything else - the actual
#     match record - moves to the other row of the pair).
#  2. Appends two brand-new match rows (239 and 240) at the bottom of the
#     sheet with upcoming fixtures that have no result yet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($rowA, $rowB) {
    $rngA = $ws.Range("B" + $rowA + ":AC" + $rowA)
    $rngB = $ws.Range("B" + $rowB + ":AC" + $rowB)
    $valA = $rngA.Value2()
    $valB = $rngB.Value2()
    $rngA.Value = $valB
    $rngB.Value = $valA
}

# Pairs of rows whose match-data needs to be swapped
Swap-Rows 123 124
Swap-Rows 180 182
Swap-Rows 187 188
Swap-Rows 189 190
Swap-Rows 195 196
Swap-Rows 199 200
Swap-Rows 205 206
Swap-Rows 211 212
Swap-Rows 217 218
Swap-Rows 230 231

# Bring formatting (styles) for the two new rows down from the last
# existing data row (238), then fill in the values.
$ws.Range("A238:AC238").Copy()
$ws.Range("A239:AC240").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 239 - Carl Zeiss Jena vs Cottbus (not yet played)
$ws.Range("A239").Value = 237
$ws.Range("B239").Value = 6880924
$ws.Range("C239").Value = "Germany Regionalliga North East"
$ws.Range("D239").Value = "Germany Regionalliga North East"
$ws.Range("E239").Value = 45387.54166666666
$ws.Range("F239").Value = "Carl Zeiss Jena"
$ws.Range("G239").Value = "Cottbus"
$ws.Range("K239").Value = 2.2
$ws.Range("L239").Value = 3.6
$ws.Range("M239").Value = 2.6
$ws.Range("N239").Value = 2.7
$ws.Range("O239").Value = 3.75
$ws.Range("P239").Value = 2.2
$ws.Range("Q239").Value = 0
$ws.Range("R239").Value = 2.15
$ws.Range("S239").Value = 1.725
$ws.Range("T239").Value = 2.75
$ws.Range("U239").Value = 1.925
$ws.Range("V239").Value = 1.925
$ws.Range("W239").Value = 0
$ws.Range("X239").Value = 0
$ws.Range("Y239").Value = 0
$ws.Range("Z239").Value = 0
$ws.Range("AA239").Value = 0
$ws.Range("H239:J239").Clear()
$ws.Range("AB239:AC239").Clear()

# Row 240 - BFC Dynamo vs Hertha Berlin II (not yet played)
$ws.Range("A240").Value = 238
$ws.Range("B240").Value = 6880504
$ws.Range("C240").Value = "Germany Regionalliga North East"
$ws.Range("D240").Value = "Germany Regionalliga North East"
$ws.Range("E240").Value = 45387.58333333334
$ws.Range("F240").Value = "BFC Dynamo"
$ws.Range("G240").Value = "Hertha Berlin II"
$ws.Range("K240").Value = 1.6
$ws.Range("L240").Value = 4
$ws.Range("M240").Value = 4.2
$ws.Range("N240").Value = 1.5
$ws.Range("O240").Value = 4.2
$ws.Range("P240").Value = 4.75
$ws.Range("Q240").Value = -1
$ws.Range("R240").Value = 1.925
$ws.Range("S240").Value = 1.925
$ws.Range("T240").Value = 3
$ws.Range("U240").Value = 1.975
$ws.Range("V240").Value = 1.875
$ws.Range("W240").Value = 0
$ws.Range("X240").Value = 0
$ws.Range("Y240").Value = 0
$ws.Range("Z240").Value = 0
$ws.Range("AA240").Value = 0
$ws.Range("H240:J240").Clear()
$ws.Range("AB240:AC240").Clear()
